# Updated cryptos list on Wed Oct 18 14:14:36 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is purely numeric-looking text (e.g. "44.45") must be
# forced to Text format first, otherwise Excel auto-converts the assigned
# string into a number (losing formatting like trailing zeros, e.g. "0.0480").
$numericLookingCells = @('D5', 'D8', 'D12', 'D19', 'D23', 'D24', 'D26', 'D27', 'D28', 'D31', 'D32', 'D39', 'D47', 'D48', 'D51')
foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).NumberFormat = '@'
}

$ws.Range('D2').Value = '28.364.68'
$ws.Range('E2').Value = '  -0.06%  '
$ws.Range('D3').Value = '1.576.87'
$ws.Range('E4').Value = '  +0.22%  '
$ws.Range('D5').Value = '212.32'
$ws.Range('E5').Value = '  +0.60%  '
$ws.Range('E6').Value = '  +0.05%  '
$ws.Range('D8').Value = '44.45'
$ws.Range('E8').Value = '  -3.37%  '
$ws.Range('E9').Value = '  +0.19%  '
$ws.Range('E10').Value = '  -0.17%  '
$ws.Range('E11').Value = '  -0.69%  '
$ws.Range('D12').Value = '0.0895'
$ws.Range('E12').Value = '  +1.84%  '
$ws.Range('D13').Value = '1.801.18'
$ws.Range('E13').Value = '  +0.14%  '
$ws.Range('D14').Value = '1.584.89'
$ws.Range('E14').Value = '  +0.77%  '
$ws.Range('E15').Value = '  -0.62%  '
$ws.Range('D16').Value = '28.379.01'
$ws.Range('E16').Value = '  -0.02%  '
$ws.Range('E17').Value = '  -1.40%  '
$ws.Range('E18').Value = '  -0.84%  '
$ws.Range('D19').Value = '231.63'
$ws.Range('E19').Value = '  +1.92%  '
$ws.Range('E20').Value = '  +0.44%  '
$ws.Range('E21').Value = '  -0.98%  '
$ws.Range('E22').Value = '  +0.09%  '
$ws.Range('D23').Value = '3.96'
$ws.Range('E23').Value = '  +0.87%  '
$ws.Range('D24').Value = '9.06'
$ws.Range('E24').Value = '  -1.42%  '
$ws.Range('E25').Value = '  +2.23%  '
$ws.Range('D26').Value = '151.62'
$ws.Range('E26').Value = '  +0.40%  '
$ws.Range('D27').Value = '15.06'
$ws.Range('E27').Value = '  +0.42%  '
$ws.Range('D28').Value = '6.39'
$ws.Range('E28').Value = '  -1.22%  '
$ws.Range('E29').Value = '  -0.51%  '
$ws.Range('E30').Value = '  +0.15%  '
$ws.Range('D31').Value = '0.0480'
$ws.Range('E31').Value = '  +3.70%  '
$ws.Range('D32').Value = '1.08'
$ws.Range('E32').Value = '  -3.71%  '
$ws.Range('E33').Value = '  -0.27%  '
$ws.Range('E34').Value = '  -1.85%  '
$ws.Range('D35').Value = '1.395.35'
$ws.Range('E35').Value = '  +0.09%  '
$ws.Range('E36').Value = '  +7.76%  '
$ws.Range('E37').Value = '  -3.30%  '
$ws.Range('E38').Value = '  +0.11%  '
$ws.Range('D39').Value = '2.66'
$ws.Range('E39').Value = '  +3.20%  '
$ws.Range('E40').Value = '  -0.70%  '
$ws.Range('E41').Value = '  -2.62%  '
$ws.Range('E42').Value = '  +0.14%  '
$ws.Range('E43').Value = '  +1.52%  '
$ws.Range('E44').Value = '  -1.14%  '
$ws.Range('E45').Value = '  -3.76%  '
$ws.Range('E46').Value = '  -2.18%  '
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').Value = '62.52'
$ws.Range('E47').Value = '  +0.23%  '
$ws.Range('B48').Value = 'WEMIXToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D48').Value = '0.926'
$ws.Range('E48').Value = '  -5.46%  '
$ws.Range('D49').Value = '1.714.39'
$ws.Range('E49').Value = '  +0.31%  '
$ws.Range('E50').Value = '  +0.10%  '
$ws.Range('B51').Value = 'BitcoinSV'
$ws.Range('C51').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D51').Value = '43.31'
$ws.Range('E51').Value = '  +10.58%  '

# Restore the cells we forced to Text back to the default (General/Normal)
# style so the saved workbook doesn't pick up stray formatting vs. the original.
foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).Style = 'Normal'
}
